$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Day 1: catering project intro - update "Project Work" and "File name/module name" notes
$ws.Range("I4").Value = "Created project base structure`nInitialized project intro file"
$ws.Range("R4").Value = "main.py`nproject_intro.py"

# Row 4 is now shorter content -> shrink the row height
$ws.Rows.Item(4).RowHeight = 150

# Update the active selection / scroll position on the sheet view
$ws.Range("R4").Select()
